$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Applicant ID" -> "Applicant code"
$ws.Range("A1").Value = "Applicant code"

# Applicant identifier values: numeric 1 -> text "NIRC6633540063e77"
$ws.Range("A2").Value = "NIRC6633540063e77"
$ws.Range("A3").Value = "NIRC6633540063e77"

# Restore the view: scroll back to show column A, select D8
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D8").Select()
